$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add manual review notes for row 7: manualAudit = 1, manualStatus = "[512]"
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = "[512]"

# Move the active cell selection to F7, matching the manual-review edit
$ws.Range("F7").Select()
